# Add a new row (row 4) to the equations sheet:
#   3^63 -> 1.144561e+30 at 1648537627856
#
# All three values must land as *text* (same as the existing rows), not as
# numbers -- even though "1.144561e+30" and "1648537627856" look numeric.
# Pre-setting NumberFormat to "@" (Text) forces Excel to store the literal
# string instead of parsing it into a number; we then copy the (unformatted)
# style from an existing data row back onto the new cells so the new row
# ends up with the same "no explicit style" look as rows 2 and 3.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").NumberFormat = "@"
$ws.Range("B4").NumberFormat = "@"
$ws.Range("C4").NumberFormat = "@"

$ws.Range("A4").Value = "3^63"
$ws.Range("B4").Value = "1.144561e+30"
$ws.Range("C4").Value = "1648537627856"

# Match the plain (unstyled) look of the other data rows.
$ws.Range("A4:C4").Style = $ws.Range("A2").Style
